$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header style (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-28
$values = @(
    @(6, 7),
    @(7, 8),
    @(6, 6),
    @(9, 9),
    @(9, 9),
    @(6, 6),
    @(6, 7),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(2, 4),
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(4, 5),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(7, 8),
    @(7, 7),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(9, 9)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
